# Rename the first column from "idCommune" to "stationId", and clean up the
# data now that the station id is looked up by name:
#   - the station id column must stay text (it's an id, not a quantity), so
#     re-enter it with a leading apostrophe the way a user would in Excel
#   - temperature columns use a comma as decimal separator (French locale)
#     instead of a dot

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: idCommune -> stationId
$ws.Range("A1").Value2 = "stationId"

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1
$firstCol = $used.Column
$lastCol = $firstCol + $used.Columns.Count - 1

for ($r = $firstRow + 1; $r -le $lastRow; $r++) {

    # Column A (station id) must be stored as text, not a number. Prefix
    # with an apostrophe (same as typing '79049004 into the cell in Excel)
    # so it is stored as text without touching the cell's number format.
    $idCell = $ws.Cells.Item($r, 1)
    $idText = $idCell.Text
    if ($idText -ne "") {
        $idCell.Value2 = "'" + $idText
    }

    # Columns D/E (temperature min/max): swap '.' for ',' as decimal mark.
    for ($c = 4; $c -le 5; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $text = $cell.Text
        if ($text -ne "" -and $text.Contains(".")) {
            $cell.Value2 = $text.Replace(".", ",")
        }
    }
}
